$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [string][char]0x2083

$updates = @(
    @(2, 4, '68.859.69', $true),
    @(2, 5, '  +0.95%  ', $true),
    @(3, 4, '3.876.51', $true),
    @(3, 5, '  -1.07%  ', $true),
    @(4, 5, '  +0.24%  ', $true),
    @(5, 4, '519.59', $true),
    @(5, 5, '  +6.72%  ', $true),
    @(6, 4, '142.33', $true),
    @(6, 5, '  -2.82%  ', $true),
    @(7, 4, '0.605', $true),
    @(7, 5, '  -2.86%  ', $true),
    @(8, 4, '0.997', $true),
    @(8, 5, '  -0.10%  ', $true),
    @(9, 4, '0.713', $true),
    @(9, 5, '  -2.77%  ', $true),
    @(10, 4, '0.168', $true),
    @(10, 5, '  +0.01%  ', $true),
    @(11, 4, '0.0000326', $true),
    @(11, 5, '  -6.49%  ', $true),
    @(12, 4, '41.61', $true),
    @(12, 5, '  -3.79%  ', $true),
    @(13, 4, '4.518.01', $true),
    @(13, 5, '  -0.59%  ', $true),
    @(14, 4, '10.11', $true),
    @(14, 5, '  -5.62%  ', $true),
    @(15, 4, '3.887.98', $true),
    @(15, 5, '  -0.51%  ', $true),
    @(16, 5, '  -0.63%  ', $true),
    @(17, 5, '  +5.90%  ', $true),
    @(18, 4, '13.64', $true),
    @(18, 5, '  -4.28%  ', $true),
    @(19, 4, '19.51', $true),
    @(19, 5, '  -2.88%  ', $true),
    @(20, 4, '68.868.66', $true),
    @(20, 5, '  +0.78%  ', $true),
    @(21, 4, '421.64', $true),
    @(21, 5, '  -2.26%  ', $true),
    @(22, 4, '3.30', $true),
    @(22, 5, '  -5.90%  ', $true),
    @(23, 4, '14.02', $true),
    @(23, 5, '  -7.47%  ', $true),
    @(24, 2, 'Litecoin', $false),
    @(24, 3, 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', $false),
    @(24, 4, '86.92', $true),
    @(24, 5, '  -1.67%  ', $true),
    @(25, 2, 'PancakeSwap', $false),
    @(25, 3, 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', $false),
    @(25, 4, '4.03', $true),
    @(25, 5, '  +8.14%  ', $true),
    @(26, 4, '11.35', $true),
    @(26, 5, '  -3.19%  ', $true),
    @(27, 4, '10.44', $true),
    @(27, 5, '  -6.67%  ', $true),
    @(28, 4, '35.58', $true),
    @(28, 5, '  -5.95%  ', $true),
    @(29, 4, '693.27', $true),
    @(29, 5, '  -3.47%  ', $true),
    @(30, 4, '12.98', $true),
    @(30, 5, '  -5.48%  ', $true),
    @(31, 4, '0.124', $true),
    @(31, 5, '  -5.26%  ', $true),
    @(32, 5, '  -4.47%  ', $true),
    @(33, 4, '67.63', $true),
    @(33, 5, '  +10.89%  ', $true),
    @(34, 4, '0.441', $true),
    @(34, 5, '  +10.97%  ', $true),
    @(35, 4, '5.87', $true),
    @(35, 5, '  -5.17%  ', $true),
    @(36, 4, '39.52', $true),
    @(36, 5, '  -5.06%  ', $true),
    @(37, 4, ('0.0' + $sub3 + '0830'), $true),
    @(37, 5, '  -9.50%  ', $true),
    @(38, 4, '0.995', $true),
    @(38, 5, '  -0.26%  ', $true),
    @(39, 5, '  +0.21%  ', $true),
    @(40, 5, '  +0.26%  ', $true),
    @(41, 4, '0.0472', $true),
    @(41, 5, '  -3.98%  ', $true),
    @(42, 4, '3.02', $true),
    @(42, 5, '  +0.00%  ', $true),
    @(43, 4, '2.72', $true),
    @(43, 5, '  -9.27%  ', $true),
    @(44, 4, '2.91', $true),
    @(44, 5, '  -6.99%  ', $true),
    @(45, 4, '3.33', $true),
    @(45, 5, '  -0.80%  ', $true),
    @(46, 2, 'EnergySwap', $false),
    @(46, 3, 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', $false),
    @(46, 4, '28.29', $true),
    @(46, 5, '  +11.71%  ', $true),
    @(47, 2, 'Stellar', $false),
    @(47, 3, 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', $false),
    @(47, 4, '0.139', $true),
    @(47, 5, '  -2.48%  ', $true),
    @(48, 2, 'Stacks', $false),
    @(48, 3, 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', $false),
    @(48, 4, '3.01', $true),
    @(48, 5, '  +7.20%  ', $true),
    @(49, 2, 'LidoDAOToken', $false),
    @(49, 3, 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', $false),
    @(49, 4, '3.25', $true),
    @(49, 5, '  -5.34%  ', $true),
    @(50, 2, 'Monero', $false),
    @(50, 3, 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', $false),
    @(50, 4, '141.53', $true),
    @(50, 5, '  -2.26%  ', $true),
    @(51, 2, 'ARBITRUM', $false),
    @(51, 3, 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', $false),
    @(51, 4, '2.02', $true),
    @(51, 5, '  -5.74%  ', $true)
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $isNumericLike = $u[3]
    $c = $ws.Cells.Item($row, $col)
    if ($isNumericLike) {
        $c.NumberFormat = "@"
        $c.Value = $val
        $c.Style = "Normal"
    } else {
        $c.Value = $val
    }
}